$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell B10: Taxonsorteringsordning 88967 -> 89101 ---
$ws.Cells.Item(10, 2).Value = 89101

# --- Add new row 11 ---
$row = 11

$ws.Cells.Item($row, 1).Value = 112365064          # A  Id
$ws.Cells.Item($row, 2).Value = 95693               # B  Taxonsorteringsordning
$ws.Cells.Item($row, 3).Value = "Ovaliderad"        # C  Valideringsstatus
$ws.Cells.Item($row, 4).Value = "LC"                # D  Rödlistade
$ws.Cells.Item($row, 5).Value = 221941              # E  TaxonId
$ws.Cells.Item($row, 6).Value = "Plattlummer"       # F  Artnamn
$ws.Cells.Item($row, 7).Value = "Lycopodium complanatum"  # G  Vetenskapligt namn
$ws.Cells.Item($row, 8).Value = "L."                # H  Auktor

# I  Antal -- present but blank in the source data
$ws.Cells.Item($row, 9).NumberFormat = "@"
$ws.Cells.Item($row, 9).Value = ""
$ws.Cells.Item($row, 9).Style = "Normal"

$ws.Cells.Item($row, 16).Value = "Flarken, Nb"       # P  Lokalnamn
$ws.Cells.Item($row, 17).Value = 813211              # Q  Ost
$ws.Cells.Item($row, 18).Value = 7316203             # R  Nord
$ws.Cells.Item($row, 19).Value = 25                  # S  Noggrannhet
$ws.Cells.Item($row, 20).Value = "Norrbotten"        # T  Län
$ws.Cells.Item($row, 21).Value = "Boden"             # U  Kommun
$ws.Cells.Item($row, 22).Value = "Norrbotten"        # V  Provins
$ws.Cells.Item($row, 23).Value = "Överluleå"         # W  Församling

# Y  Startdatum -- stored as plain text, not a date serial
$ws.Cells.Item($row, 25).NumberFormat = "@"
$ws.Cells.Item($row, 25).Value = "2023-09-05"
$ws.Cells.Item($row, 25).Style = "Normal"

# AA Slutdatum -- stored as plain text, not a date serial
$ws.Cells.Item($row, 27).NumberFormat = "@"
$ws.Cells.Item($row, 27).Value = "2023-09-05"
$ws.Cells.Item($row, 27).Style = "Normal"

$ws.Cells.Item($row, 30).Value = $false              # AD Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false              # AE Osäker artbestämning
$ws.Cells.Item($row, 33).Value = $false              # AG Ospontan

# AT Bestämningsår -- present but blank in the source data
$ws.Cells.Item($row, 46).NumberFormat = "@"
$ws.Cells.Item($row, 46).Value = ""
$ws.Cells.Item($row, 46).Style = "Normal"

$ws.Cells.Item($row, 49).Value = "Linnea Åsedahl"    # AW Rapportör
$ws.Cells.Item($row, 50).Value = "Linnea Åsedahl"    # AX Observatörer

# AY Projektnamn -- present but blank in the source data
$ws.Cells.Item($row, 51).NumberFormat = "@"
$ws.Cells.Item($row, 51).Value = ""
$ws.Cells.Item($row, 51).Style = "Normal"

Write-Output "Row 11 added and B10 updated"
